$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.166.36'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '2.930.39'
$ws.Range("E3").Value = '  +1.19%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.87'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.49'
$ws.Range("E6").Value = '  +1.40%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.506'
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.97'
$ws.Range("E9").Value = '  +4.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.441'
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.77'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '3.415.54'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '61.138.73'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.73'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").Value = '2.933.82'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '436.82'
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.58'
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("E24").Value = '  +3.01%  '
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.91'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  +3.13%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.02'
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("E31").Value = '  +4.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.69'
$ws.Range("E32").Value = '  +1.61%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").Value = '0.0₃0870'
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.00'
$ws.Range("E38").Value = '  +1.29%  '
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '42.16'
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '376.85'
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0348'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = '2.691.96'
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.52'
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.07'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.01'
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  +1.12%  '
